$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains text formatting so numeric-looking
# values like "243.17" or "0.05985" are stored as text, not converted to numbers.
$ws.Range("D2:D48").NumberFormat = "@"

$ws.Range("D2").Value = '243.17'
$ws.Range("D3").Value = '23.12'
$ws.Range("D4").Value = '5.406'
$ws.Range("D5").Value = '0.05985'
$ws.Range("D6").Value = '3.430'
$ws.Range("D7").Value = '6.519'
$ws.Range("D8").Value = '0.8083'
$ws.Range("D9").Value = '0.9211'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1428'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07431'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03289'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03068'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09357'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '3.846'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001576'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04710'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = '0.0005921'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("D19").Value = '0.005857'
$ws.Range("D21").Value = '0.004875'
$ws.Range("D23").Value = '3.569'
$ws.Range("D24").Value = '2.135'
$ws.Range("E27").Value = '26UpBotsUBXT'
$ws.Range("D40").Value = '0.03971'
$ws.Range("B41").Value = 'CEJI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D41").Value = '0.004901'
$ws.Range("E41").Value = '40CEJICEJIBestin24h'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '0.006441'
$ws.Range("E42").Value = '41KickTokenKICK'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = '0.1076'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").Value = '0.009204'
$ws.Range("D45").Value = '0.00005080'
$ws.Range("D47").Value = '0.7002'
$ws.Range("D48").Value = '0.002433'
